$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataCell($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 10
    $rng.Font.ThemeColor = 1
    if ($null -ne $text) {
        $rng.Value = $text
    }
}

# Row 2 - MCH197-1 series
Set-DataCell "A2" "MCH197-1"
Set-DataCell "C2" "FACT SHEETS ON SWEEDEN, SWEDISH PARLIAMENTARY SYSTEM- MARKET POLICY, LOCAL GOVERNMENT IN SWEDEN"
Set-DataCell "D2" $null
Set-DataCell "E2" "Series"
Set-DataCell "F2" "1 Box"
Set-DataCell "G2" "LOCATION: 24B | GRAP COUNT NUMER: NONE"
Set-DataCell "H2" $null

# Row 3 - MCH197-2 series
Set-DataCell "A3" "MCH197-2"
Set-DataCell "C3" "GENERAL"
Set-DataCell "D3" $null
Set-DataCell "E3" "Series"
Set-DataCell "F3" "1 Box"
Set-DataCell "G3" "LOCATION: 24B | GRAP COUNT NUMER: NONE"
Set-DataCell "H3" $null

# Restore frozen header pane (row 1 frozen) and selection over the new data rows
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A2:M3").Select() | Out-Null
